# Auto-generated edit script: updates market-price derived columns (H-N)
# across the ALC/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled runner's
# refreshed pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 343
$ws.Range("I48").Value = 217
$ws.Range("K48").Value = 651
$ws.Range("M48").Value = -359
$ws.Range("H56").Value = 343
$ws.Range("I56").Value = 217
$ws.Range("K56").Value = 651
$ws.Range("M56").Value = -117
$ws.Range("H63").Value = 50000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 50000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 50000
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -51248
$ws.Range("H66").Value = 50000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 50000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 150000
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -156240
$ws.Range("H92").Value = 1409.375
$ws.Range("I92").Value = 1490.0667
$ws.Range("K92").Value = 1490.0667
$ws.Range("M92").Value = -242.0667000000001
$ws.Range("H96").Value = 1890.5217
$ws.Range("I96").Value = 887.6
$ws.Range("J96").Value = 2662
$ws.Range("K96").Value = 2662.8
$ws.Range("L96").Value = 7986
$ws.Range("M96").Value = -1289.8
$ws.Range("N96").Value = -10732
$ws.Range("H97").Value = 50109.5
$ws.Range("J97").Value = 50109.5
$ws.Range("L97").Value = 150328.5
$ws.Range("N97").Value = -151320.5
$ws.Range("H98").Value = 960.0454999999999
$ws.Range("I98").Value = 786.41174
$ws.Range("K98").Value = 786.41174
$ws.Range("M98").Value = 711.58826
$ws.Range("H100").Value = 3598.8857
$ws.Range("J100").Value = 5286.1577
$ws.Range("L100").Value = 5286.1577
$ws.Range("N100").Value = -6368.1577
$ws.Range("H122").Value = 960.0454999999999
$ws.Range("I122").Value = 786.41174
$ws.Range("K122").Value = 2359.23522
$ws.Range("M122").Value = 90.76477999999997
$ws.Range("H138").Value = 2499.7307
$ws.Range("J138").Value = 3200
$ws.Range("L138").Value = 9600
$ws.Range("N138").Value = -19880

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 65000
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H94").Value = 6574.6924
$ws.Range("I94").Value = 6988.5557
$ws.Range("J94").Value = 5643.5
$ws.Range("K94").Value = 6988.5557
$ws.Range("L94").Value = 5643.5
$ws.Range("M94").Value = -6537.5557
$ws.Range("N94").Value = -6545.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6888.737
$ws.Range("I58").Value = 5818.625
$ws.Range("K58").Value = 5818.625
$ws.Range("M58").Value = -5615.625
$ws.Range("H136").Value = 6888.737
$ws.Range("I136").Value = 5818.625
$ws.Range("K136").Value = 17455.875
$ws.Range("M136").Value = -14905.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 717451.9
$ws.Range("J32").Value = 913031.9399999999
$ws.Range("L32").Value = 2739095.82
$ws.Range("N32").Value = -2739661.82
$ws.Range("H105").Value = 16200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 4125.8335
$ws.Range("I21").Value = 4125.8335
$ws.Range("K21").Value = 4125.8335
$ws.Range("M21").Value = -3952.8335
$ws.Range("H30").Value = 4125.8335
$ws.Range("I30").Value = 4125.8335
$ws.Range("K30").Value = 4125.8335
$ws.Range("M30").Value = -4020.8335
$ws.Range("H70").Value = 14264.842
$ws.Range("I70").Value = 15268.8
$ws.Range("K70").Value = 15268.8
$ws.Range("M70").Value = -14998.8
$ws.Range("H73").Value = 14264.842
$ws.Range("I73").Value = 15268.8
$ws.Range("K73").Value = 15268.8
$ws.Range("M73").Value = -14332.8
$ws.Range("H126").Value = 2968.7222
$ws.Range("I126").Value = 3143.7
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 9431.099999999999
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = -6961.099999999999
$ws.Range("N126").Value = -13190

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2782.0244
$ws.Range("I46").Value = 1384.6154
$ws.Range("J46").Value = 3430.8215
$ws.Range("K46").Value = 1384.6154
$ws.Range("L46").Value = 3430.8215
$ws.Range("M46").Value = -1196.6154
$ws.Range("N46").Value = -3806.8215
$ws.Range("H61").Value = 6840.4814
$ws.Range("I61").Value = 7581.8945
$ws.Range("J61").Value = 5079.625
$ws.Range("K61").Value = 7581.8945
$ws.Range("L61").Value = 5079.625
$ws.Range("M61").Value = -7379.8945
$ws.Range("N61").Value = -5483.625
$ws.Range("H93").Value = 8391.5
$ws.Range("I93").Value = 7032.6
$ws.Range("K93").Value = 7032.6
$ws.Range("M93").Value = -5784.6
$ws.Range("H113").Value = 6840.4814
$ws.Range("I113").Value = 7581.8945
$ws.Range("J113").Value = 5079.625
$ws.Range("K113").Value = 7581.8945
$ws.Range("L113").Value = 5079.625
$ws.Range("M113").Value = -5411.8945
$ws.Range("N113").Value = -9419.625
$ws.Range("H133").Value = 88777
$ws.Range("J133").Value = 88777
$ws.Range("L133").Value = 88777
$ws.Range("N133").Value = -93837

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 17571.428
$ws.Range("I2").Value = 17571.428
$ws.Range("K2").Value = 17571.428
$ws.Range("M2").Value = -17459.428
$ws.Range("H4").Value = 105699.8
$ws.Range("I4").Value = 116666.11
$ws.Range("K4").Value = 116666.11
$ws.Range("M4").Value = -116553.11
$ws.Range("H8").Value = 1502500
$ws.Range("J8").Value = 5000
$ws.Range("L8").Value = 5000
$ws.Range("N8").Value = -5280
$ws.Range("H14").Value = 3000
$ws.Range("I14").Value = 3000
$ws.Range("K14").Value = 3000
$ws.Range("M14").Value = -2832
$ws.Range("H18").Value = 35000
$ws.Range("I18").Value = 30000
$ws.Range("J18").Value = 40000
$ws.Range("K18").Value = 30000
$ws.Range("L18").Value = 40000
$ws.Range("M18").Value = -29827
$ws.Range("N18").Value = -40346
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()
$ws.Range("H100").Value = 1164.9166
$ws.Range("I100").Value = 907.9
$ws.Range("K100").Value = 1815.8
$ws.Range("M100").Value = -1274.8

